# Update cryptocurrency price/volume figures per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.818.93'
$ws.Range("E2").Value = '  +0.43%  '
$ws.Range("D3").Value = '1.905.90'
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'312.95"
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").Value = "'0.5213"
$ws.Range("E7").Value = '  +7.18%  '
$ws.Range("D8").Value = "'0.3790"
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").Value = "'0.07238"
$ws.Range("E9").Value = '  -1.16%  '
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = "'0.9114"
$ws.Range("E10").Value = '  -0.25%  '
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").Value = "'21.26"
$ws.Range("E11").Value = '  +3.49%  '
$ws.Range("D12").Value = '1.933.37'
$ws.Range("E12").Value = '  +2.37%  '
$ws.Range("D13").Value = "'0.07647"
$ws.Range("E13").Value = '  -0.18%  '
$ws.Range("D14").Value = "'5.451"
$ws.Range("E14").Value = '  -0.51%  '
$ws.Range("D15").Value = "'92.20"
$ws.Range("E15").Value = '  +1.15%  '
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = "'0.000008703"
$ws.Range("E17").Value = '  -0.84%  '
$ws.Range("D18").Value = "'1.000"
$ws.Range("E18").Value = '  -0.12%  '
$ws.Range("D19").Value = '27.859.66'
$ws.Range("E19").Value = '  +1.04%  '
$ws.Range("D20").Value = "'14.54"
$ws.Range("E20").Value = '  +0.32%  '
$ws.Range("D21").Value = "'5.153"
$ws.Range("E21").Value = '  +0.66%  '
$ws.Range("D22").Value = '2.172.39'
$ws.Range("E22").Value = '  +1.90%  '
$ws.Range("D23").Value = "'10.86"
$ws.Range("E23").Value = '  +1.11%  '
$ws.Range("D24").Value = "'6.641"
$ws.Range("E24").Value = '  +0.68%  '
$ws.Range("D25").Value = "'153.55"
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("D26").Value = "'1.866"
$ws.Range("E26").Value = '  -2.27%  '
$ws.Range("D27").Value = "'2.168"
$ws.Range("E27").Value = '  +0.87%  '
$ws.Range("D28").Value = "'18.33"
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("D29").Value = "'114.74"
$ws.Range("E29").Value = '  -0.62%  '
$ws.Range("D30").Value = "'4.860"
$ws.Range("E30").Value = '  -0.12%  '
$ws.Range("D31").Value = "'0.09008"
$ws.Range("E31").Value = '  +1.32%  '
$ws.Range("D32").Value = "'4.865"
$ws.Range("E32").Value = '  +5.11%  '
$ws.Range("D33").Value = "'3.181"
$ws.Range("E33").Value = '  -0.70%  '
$ws.Range("E34").Value = '  +1.11%  '
$ws.Range("D35").Value = "'0.7800"
$ws.Range("E35").Value = '  +1.88%  '
$ws.Range("D36").Value = "'0.02097"
$ws.Range("E36").Value = '  +2.92%  '
$ws.Range("D37").Value = "'2.606"
$ws.Range("E37").Value = '  +3.24%  '
$ws.Range("D38").Value = "'3.075"
$ws.Range("E38").Value = '  +3.34%  '
$ws.Range("D39").Value = "'0.5577"
$ws.Range("E39").Value = '  +1.95%  '
$ws.Range("D40").Value = "'1.093"
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("D41").Value = "'0.05284"
$ws.Range("E41").Value = '  +0.25%  '
$ws.Range("D42").Value = "'6.719"
$ws.Range("E42").Value = '  -2.44%  '
$ws.Range("D43").Value = "'115.71"
$ws.Range("E43").Value = '  +3.57%  '
$ws.Range("D44").Value = "'8.554"
$ws.Range("E44").Value = '  +0.66%  '
$ws.Range("D45").Value = "'0.1516"
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("D46").Value = "'0.4813"
$ws.Range("D47").Value = "'10.46"
$ws.Range("E47").Value = '  -1.58%  '
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("D49").Value = "'1.619"
$ws.Range("E49").Value = '  -0.78%  '
$ws.Range("D50").Value = "'67.00"
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("D51").Value = "'0.05995"
$ws.Range("E51").Value = '  -0.90%  '
